# Insert a new weekly price record at the top of the "Feria Lagunitas de
# Puerto Montt - Zanahoria" data block (row 544), pushing the existing
# rows 544:612 down to 545:613.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(544).Insert()

$ws.Cells.Item(544, 1).Value  = 4
$ws.Cells.Item(544, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(544, 3).Value  = "Los Lagos"
$ws.Cells.Item(544, 4).Value  = 45142
$ws.Cells.Item(544, 5).Value  = 10
$ws.Cells.Item(544, 6).Value  = 100114013
$ws.Cells.Item(544, 7).Value  = "Zanahoria"
$ws.Cells.Item(544, 8).Value  = "Sin especificar"
$ws.Cells.Item(544, 9).Value  = "Primera"
$ws.Cells.Item(544, 10).Value = 900
$ws.Cells.Item(544, 11).Value = 7500
$ws.Cells.Item(544, 12).Value = 9000
$ws.Cells.Item(544, 13).Value = 8250
$ws.Cells.Item(544, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(544, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(544, 16).Value = 412
$ws.Cells.Item(544, 17).Value = 20
$ws.Cells.Item(544, 18).Value = "Hortaliza"
